$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value would otherwise be re-interpreted as a
# number by Excel's smart-typing (e.g. "610.80" -> 610.8); force text first
# so the stored value matches the source data exactly.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.989.66'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '3.154.42'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '610.80'
$ws.Range("E5").Value = '  +2.36%  '
$ws.Range("D6").Value = '146.99'
$ws.Range("E6").Value = '  -3.28%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.147.26'
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("E10").Value = '  -1.59%  '
$ws.Range("D11").Value = '5.39'
$ws.Range("E11").Value = '  -3.23%  '
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  -1.29%  '
$ws.Range("D14").Value = '35.56'
$ws.Range("E14").Value = '  -3.66%  '
$ws.Range("D15").Value = '3.669.30'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("E16").Value = '  +3.04%  '
$ws.Range("D17").Value = '63.976.76'
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").Value = '3.150.50'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").Value = '6.91'
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("D20").Value = '477.47'
$ws.Range("E20").Value = '  -1.24%  '
$ws.Range("D21").Value = '14.64'
$ws.Range("E21").Value = '  -0.90%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '7.99'
$ws.Range("E23").Value = '  +2.89%  '
$ws.Range("D24").Value = '13.79'
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("E25").Value = '  -1.89%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  -3.45%  '
$ws.Range("D28").Value = '8.53'
$ws.Range("E28").Value = '  -1.70%  '
$ws.Range("E29").Value = '  +2.54%  '
$ws.Range("D30").Value = '0.120'
$ws.Range("E30").Value = '  -6.02%  '
$ws.Range("D31").Value = '2.11'
$ws.Range("E31").Value = '  -7.31%  '
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("D34").Value = '26.32'
$ws.Range("E34").Value = '  -2.14%  '
$ws.Range("E35").Value = '  +1.83%  '
$ws.Range("D36").Value = '0.0₃0789'
$ws.Range("E36").Value = '  +8.21%  '
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("D38").Value = '52.91'
$ws.Range("E38").Value = '  -2.94%  '
$ws.Range("D39").Value = '460.06'
$ws.Range("E39").Value = '  -2.69%  '
$ws.Range("D40").Value = '3.04'
$ws.Range("E40").Value = '  -6.94%  '
$ws.Range("D41").Value = '0.0398'
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("E42").Value = '  -4.42%  '
$ws.Range("E43").Value = '  -2.11%  '
$ws.Range("D44").Value = '2.867.69'
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("E45").Value = '  -4.62%  '
$ws.Range("E46").Value = '  -2.87%  '
$ws.Range("E47").Value = '  +2.45%  '
$ws.Range("D48").Value = '26.48'
$ws.Range("E48").Value = '  -3.32%  '
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D51").Value = '118.96'
$ws.Range("E51").Value = '  -1.56%  '
